$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Pin on Dev" header labels (col E / J mirror the existing "Pin" header) ---
$ws.Range("E1").Value = "Pin on Dev"
$ws.Range("J1").Value = "Pin on Dev"

# --- Row 2: ultrasonic distance sensor +5V feed, wiring colour note ---
# (leading "'" keeps this an explicit-text edit-in-place so the cell's
#  existing quote-prefixed/centered number format (s=4) carries over)
$ws.Range("I2").Value = "'Distance +5"
$ws.Range("J2").Value = "6 Green"

# --- Row 3: pHat +5V rail (was a bare "+5" before) ---
$ws.Range("I3").Value = "'pHat + 5"

# --- Row 4: pHat GND (was a bare "GND" before) ---
$ws.Range("I4").Value = "pHat GND"

# --- Row 10: new blank-but-centered cell M10 ---
$ws.Range("M10").HorizontalAlignment = -4108

# --- Row 15: newly used "Not Usable" pins ---
$ws.Range("F15").Value = "Not Usable"
$ws.Range("F15").HorizontalAlignment = -4108
$ws.Range("I15").Value = "Not Usable"
$ws.Range("I15").HorizontalAlignment = -4108

# --- Row 16: Distance sensor (right) trig/echo wiring + GND rail rename ---
$ws.Range("E16").Value = "9 Red"
$ws.Range("F16").Value = "Distance Right Trig"
$ws.Range("F16").HorizontalAlignment = -4108
$ws.Range("I16").Value = "Distance GND"
$ws.Range("J16").Value = "3 Gray"

# --- Row 17: Distance sensor (center) echo wiring + right echo label ---
$ws.Range("E17").Value = "4 Purple"
$ws.Range("F17").Value = "Distance Center Echo"
$ws.Range("F17").HorizontalAlignment = -4108
$ws.Range("I17").Value = "Distance Right Echo"
$ws.Range("I17").HorizontalAlignment = -4108
$ws.Range("J17").Value = "10 Brown"

# --- Row 18: Distance sensor (center) trig wiring + right motor GND rename ---
$ws.Range("E18").Value = "5 Blue"
$ws.Range("F18").Value = "Distance Center Trig"
$ws.Range("F18").HorizontalAlignment = -4108
$ws.Range("I18").Value = "Right Motor GND"

# --- Row 19: Distance sensor (left) trig wiring ---
$ws.Range("E19").Value = "2 White"
$ws.Range("F19").Value = "Distance Left Trig"
$ws.Range("F19").HorizontalAlignment = -4108

# --- Row 20: Distance sensor (left) echo wiring ---
$ws.Range("E20").Value = "1 Black"
$ws.Range("F20").Value = "Distance Left Echo"
$ws.Range("F20").HorizontalAlignment = -4108

# --- Row 21: left motor GND rename (was a bare "GND" before) ---
$ws.Range("F21").Value = "Left Motor GND"

# --- Column E widened to fit the new "Distance Right Trig"-style labels ---
$ws.Columns("E").ColumnWidth = 8.74

# --- Move the active selection to reflect where editing finished ---
[void]$ws.Range("M19").Select()
